# Weekly update: a new daily price record was inserted as the new row 73
# (pushing the previously-existing rows 73..170 down to 74..171), adding a
# fresh "Zapallo italiano" observation for Vega Monumental Concepción.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 73; everything below shifts down one row.
$ws.Rows.Item(73).Insert()

# Populate the new row 73 with the new observation's data.
$ws.Cells.Item(73, 1).Value = 11
$ws.Cells.Item(73, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(73, 3).Value = "Bíobío"
$ws.Cells.Item(73, 4).Value = 44894
$ws.Cells.Item(73, 5).Value = 8
$ws.Cells.Item(73, 6).Value = 100112032
$ws.Cells.Item(73, 7).Value = "Zapallo italiano"
$ws.Cells.Item(73, 8).Value = "Sin especificar"
$ws.Cells.Item(73, 9).Value = "Primera"
$ws.Cells.Item(73, 10).Value = 100
$ws.Cells.Item(73, 11).Value = 7000
$ws.Cells.Item(73, 12).Value = 7500
$ws.Cells.Item(73, 13).Value = 7250
$ws.Cells.Item(73, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(73, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(73, 16).Value = 145
$ws.Cells.Item(73, 17).Value = 50
$ws.Cells.Item(73, 18).Value = "Hortaliza"
